$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '37.782.52'
Set-TextValue $ws.Range("E2") '  -0.14%  '
Set-TextValue $ws.Range("D3") '2.075.59'
Set-TextValue $ws.Range("E3") '  -0.57%  '
Set-TextValue $ws.Range("E4") '  +0.12%  '
Set-TextValue $ws.Range("D5") '233.41'
Set-TextValue $ws.Range("E5") '  -0.69%  '
Set-TextValue $ws.Range("E6") '  -0.27%  '
Set-TextValue $ws.Range("D7") '58.44'
Set-TextValue $ws.Range("E7") '  -1.85%  '
Set-TextValue $ws.Range("E8") '  +0.02%  '
Set-TextValue $ws.Range("E9") '  +0.10%  '
Set-TextValue $ws.Range("D10") '0.0785'
Set-TextValue $ws.Range("E10") '  -0.88%  '
Set-TextValue $ws.Range("E11") '  +3.71%  '
Set-TextValue $ws.Range("D12") '2.381.77'
Set-TextValue $ws.Range("E12") '  -0.52%  '
Set-TextValue $ws.Range("D13") '14.75'
Set-TextValue $ws.Range("E13") '  +0.21%  '
Set-TextValue $ws.Range("D14") '21.11'
Set-TextValue $ws.Range("E14") '  -1.47%  '
Set-TextValue $ws.Range("E15") '  +0.57%  '
Set-TextValue $ws.Range("E16") '  -0.02%  '
Set-TextValue $ws.Range("D17") '2.073.01'
Set-TextValue $ws.Range("E17") '  -0.69%  '
Set-TextValue $ws.Range("D18") '37.688.82'
Set-TextValue $ws.Range("E18") '  -0.17%  '
Set-TextValue $ws.Range("D19") '6.14'
Set-TextValue $ws.Range("E19") '  -1.16%  '
Set-TextValue $ws.Range("D20") '71.50'
Set-TextValue $ws.Range("E20") '  -0.24%  '
Set-TextValue $ws.Range("D21") '0.0₃0842'
Set-TextValue $ws.Range("E21") '  +1.42%  '
Set-TextValue $ws.Range("D22") '229.08'
Set-TextValue $ws.Range("E22") '  +0.06%  '
Set-TextValue $ws.Range("D23") '1.00'
Set-TextValue $ws.Range("E23") '  -0.02%  '
Set-TextValue $ws.Range("D24") '2.39'
Set-TextValue $ws.Range("E24") '  -0.83%  '
Set-TextValue $ws.Range("E25") '  -0.49%  '
Set-TextValue $ws.Range("D26") '9.71'
Set-TextValue $ws.Range("E26") '  +7.16%  '
Set-TextValue $ws.Range("D27") '171.79'
Set-TextValue $ws.Range("E27") '  +0.77%  '
Set-TextValue $ws.Range("E28") '  -0.18%  '
Set-TextValue $ws.Range("B29") 'EthereumClassic'
Set-TextValue $ws.Range("C29") 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range("D29") '19.41'
Set-TextValue $ws.Range("E29") '  -0.92%  '
Set-TextValue $ws.Range("B30") 'ImmutableX'
Set-TextValue $ws.Range("C30") 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range("D30") '1.40'
Set-TextValue $ws.Range("E30") '  -2.59%  '
Set-TextValue $ws.Range("E31") '  +1.25%  '
Set-TextValue $ws.Range("E32") '  -0.01%  '
Set-TextValue $ws.Range("E33") '  +0.23%  '
Set-TextValue $ws.Range("D34") '4.70'
Set-TextValue $ws.Range("E34") '  -0.11%  '
Set-TextValue $ws.Range("E35") '  -3.05%  '
Set-TextValue $ws.Range("E36") '  -0.22%  '
Set-TextValue $ws.Range("E37") '  -3.79%  '
Set-TextValue $ws.Range("D38") '1.00'
Set-TextValue $ws.Range("E38") '  +0.22%  '
Set-TextValue $ws.Range("E39") '  -0.29%  '
Set-TextValue $ws.Range("D40") '0.0234'
Set-TextValue $ws.Range("E40") '  +8.89%  '
Set-TextValue $ws.Range("D41") '100.28'
Set-TextValue $ws.Range("E41") '  +0.08%  '
Set-TextValue $ws.Range("D42") '0.0978'
Set-TextValue $ws.Range("E42") '  -0.91%  '
Set-TextValue $ws.Range("D43") '17.19'
Set-TextValue $ws.Range("E43") '  +5.98%  '
Set-TextValue $ws.Range("D44") '2.91'
Set-TextValue $ws.Range("E44") '  -1.11%  '
Set-TextValue $ws.Range("D45") '1.448.00'
Set-TextValue $ws.Range("E45") '  -1.11%  '
Set-TextValue $ws.Range("D46") '1.16'
Set-TextValue $ws.Range("E46") '  -1.58%  '
Set-TextValue $ws.Range("E47") '  -0.46%  '
Set-TextValue $ws.Range("D48") '4.10'
Set-TextValue $ws.Range("E48") '  -2.69%  '
Set-TextValue $ws.Range("E50") '  -1.65%  '
Set-TextValue $ws.Range("D51") '2.267.39'
Set-TextValue $ws.Range("E51") '  -0.54%  '
